$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $val) {
    $c = $ws.Range($cellref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.NumberFormat = "General"
}

Set-TextValue "D2" "29.299.42"
Set-TextValue "E2" "  +0.27%  "
Set-TextValue "D3" "1.860.05"
Set-TextValue "E3" "  +0.12%  "
Set-TextValue "D4" "1.0000"
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "0.7048"
Set-TextValue "E5" "  +0.90%  "
Set-TextValue "D6" "238.27"
Set-TextValue "E6" "  +0.42%  "
Set-TextValue "D7" "0.9999"
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "0.07918"
Set-TextValue "E8" "  +2.49%  "
Set-TextValue "D9" "0.3048"
Set-TextValue "E9" "  +0.08%  "
Set-TextValue "D10" "24.58"
Set-TextValue "E10" "  +5.69%  "
Set-TextValue "B11" "WrappedEther"
Set-TextValue "C11" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D11" "2.601.21"
Set-TextValue "E11" "  +39.46%  "
Set-TextValue "B12" "TRON"
Set-TextValue "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.08182"
Set-TextValue "E12" "  +0.23%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.229"
Set-TextValue "E13" "  +1.47%  "
Set-TextValue "B14" "Polygon"
Set-TextValue "C14" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D14" "0.7192"
Set-TextValue "E14" "  +0.33%  "
Set-TextValue "B15" "Litecoin"
Set-TextValue "C15" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D15" "89.64"
Set-TextValue "E15" "  +0.57%  "
Set-TextValue "D16" "28.506.73"
Set-TextValue "E16" "  -2.43%  "
Set-TextValue "D17" "5.825"
Set-TextValue "E17" "  +1.28%  "
Set-TextValue "D18" "0.000007805"
Set-TextValue "E18" "  +0.95%  "
Set-TextValue "D19" "13.23"
Set-TextValue "E19" "  -0.61%  "
Set-TextValue "D20" "238.37"
Set-TextValue "E20" "  +0.24%  "
Set-TextValue "D21" "1.000"
Set-TextValue "E21" "  +0.06%  "
Set-TextValue "D22" "1.0000"
Set-TextValue "E22" "  +0.05%  "
Set-TextValue "D23" "7.516"
Set-TextValue "E23" "  +1.24%  "
Set-TextValue "D24" "1.939.26"
Set-TextValue "E24" "  -8.08%  "
Set-TextValue "D25" "162.70"
Set-TextValue "E25" "  +0.05%  "
Set-TextValue "D26" "8.903"
Set-TextValue "E26" "  -1.13%  "
Set-TextValue "D27" "0.1432"
Set-TextValue "E27" "  -3.47%  "
Set-TextValue "E28" "  +0.48%  "
Set-TextValue "D29" "1.923"
Set-TextValue "E29" "  -5.93%  "
Set-TextValue "D30" "1.379"
Set-TextValue "E30" "  -2.64%  "
Set-TextValue "D31" "1.474"
Set-TextValue "E31" "  -0.46%  "
Set-TextValue "D32" "4.336"
Set-TextValue "D33" "4.059"
Set-TextValue "E33" "  +0.93%  "
Set-TextValue "D34" "0.05183"
Set-TextValue "E34" "  -0.11%  "
Set-TextValue "D35" "1.176"
Set-TextValue "E35" "  +0.89%  "
Set-TextValue "D36" "0.7117"
Set-TextValue "E36" "  +0.57%  "
Set-TextValue "D37" "0.9903"
Set-TextValue "E37" "  -0.86%  "
Set-TextValue "D38" "2.674"
Set-TextValue "E38" "  +0.68%  "
Set-TextValue "E39" "  +0.30%  "
Set-TextValue "D40" "2.690"
Set-TextValue "E40" "  -1.19%  "
Set-TextValue "D41" "1.172.00"
Set-TextValue "E41" "  +2.89%  "
Set-TextValue "D42" "0.9243"
Set-TextValue "E42" "  -1.49%  "
Set-TextValue "D43" "5.961"
Set-TextValue "E43" "  +1.23%  "
Set-TextValue "D44" "0.4258"
Set-TextValue "E44" "  -0.38%  "
Set-TextValue "D45" "70.85"
Set-TextValue "E45" "  +0.27%  "
Set-TextValue "D46" "0.9995"
Set-TextValue "E46" "  -0.05%  "
Set-TextValue "D47" "101.36"
Set-TextValue "E47" "  -1.71%  "
Set-TextValue "D48" "0.5326"
Set-TextValue "E48" "  -2.60%  "
Set-TextValue "D49" "1.757"
Set-TextValue "E49" "  -1.84%  "
Set-TextValue "D50" "9.192"
Set-TextValue "E50" "  +0.36%  "
Set-TextValue "D51" "7.002"
Set-TextValue "E51" "  +0.75%  "
